$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CALENDAR")

$errorMsg = "Could not find Chrome (ver. 145.0.7632.67). This can occur if either`n 1. you did not perform an installation before running the script (e.g. ``npx puppeteer browsers install chrome``) or`n 2. your cache path is incorrectly configured (which is: /home/jules/.cache/puppeteer).`nFor (2), check out our guide on configuring puppeteer at https://pptr.dev/guides/configuration."

$ws.Range("I1").Value = "error_log"

for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Value = "failed"
    $ws.Cells.Item($r, 9).Value = $errorMsg
}
